$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 262.5
$ws.Range("I39").Value = 250
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 750
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = -454
$ws.Range("N39").Value = -1492

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 38463532
$ws.Range("I40").Value = 1250
$ws.Range("J40").Value = 45456670
$ws.Range("K40").Value = 1250
$ws.Range("L40").Value = 45456670
$ws.Range("M40").Value = -1075
$ws.Range("N40").Value = -45457020

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2072.7856
$ws.Range("I127").Value = 582.6667
$ws.Range("J127").Value = 2900.6296
$ws.Range("K127").Value = 1748.0001
$ws.Range("L127").Value = 8701.888800000001
$ws.Range("M127").Value = 3211.9999
$ws.Range("N127").Value = -18621.8888

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 993.2954999999999
$ws.Range("I135").Value = 545.37836
$ws.Range("J135").Value = 3360.8572
$ws.Range("K135").Value = 4908.40524
$ws.Range("L135").Value = 30247.7148
$ws.Range("M135").Value = -2373.40524
$ws.Range("N135").Value = -35317.7148

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 825.25
$ws.Range("I137").Value = 777.86206
$ws.Range("J137").Value = 950.1818
$ws.Range("K137").Value = 2333.58618
$ws.Range("L137").Value = 2850.5454
$ws.Range("M137").Value = 216.4138199999998
$ws.Range("N137").Value = -7950.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1748.38
$ws.Range("I141").Value = 584.619
$ws.Range("J141").Value = 7858.125
$ws.Range("K141").Value = 1753.857
$ws.Range("L141").Value = 23574.375
$ws.Range("M141").Value = 3426.143
$ws.Range("N141").Value = -33934.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17400.945
$ws.Range("I32").Value = 19809.508
$ws.Range("J32").Value = 9627.862999999999
$ws.Range("K32").Value = 19809.508
$ws.Range("L32").Value = 9627.862999999999
$ws.Range("M32").Value = -19522.508
$ws.Range("N32").Value = -10201.863

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1035.591
$ws.Range("I61").Value = 730.3158
$ws.Range("K61").Value = 730.3158
$ws.Range("M61").Value = -518.3158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 806.9048
$ws.Range("I74").Value = 726.5
$ws.Range("J74").Value = 1289.3334
$ws.Range("K74").Value = 726.5
$ws.Range("L74").Value = 1289.3334
$ws.Range("M74").Value = 147.5
$ws.Range("N74").Value = -3037.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 806.9048
$ws.Range("I77").Value = 726.5
$ws.Range("J77").Value = 1289.3334
$ws.Range("K77").Value = 3632.5
$ws.Range("L77").Value = 6446.666999999999
$ws.Range("M77").Value = 735.5
$ws.Range("N77").Value = -15182.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1253.9636
$ws.Range("I132").Value = 929.6667
$ws.Range("J132").Value = 2301.6924
$ws.Range("K132").Value = 2789.0001
$ws.Range("L132").Value = 6905.0772
$ws.Range("M132").Value = -259.0001000000002
$ws.Range("N132").Value = -11965.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1035.591
$ws.Range("I136").Value = 730.3158
$ws.Range("K136").Value = 2190.9474
$ws.Range("M136").Value = 359.0526

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19029.932
$ws.Range("I134").Value = 1419.5217
$ws.Range("J134").Value = 86536.5
$ws.Range("K134").Value = 4258.5651
$ws.Range("L134").Value = 259609.5
$ws.Range("M134").Value = -1723.5651
$ws.Range("N134").Value = -264679.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2430.1296
$ws.Range("I31").Value = 2470.2222
$ws.Range("J31").Value = 2349.9443
$ws.Range("K31").Value = 2470.2222
$ws.Range("L31").Value = 2349.9443
$ws.Range("M31").Value = -2175.2222
$ws.Range("N31").Value = -2939.9443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2430.1296
$ws.Range("I34").Value = 2470.2222
$ws.Range("J34").Value = 2349.9443
$ws.Range("K34").Value = 2470.2222
$ws.Range("L34").Value = 2349.9443
$ws.Range("M34").Value = -2268.2222
$ws.Range("N34").Value = -2753.9443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2031.5161
$ws.Range("I132").Value = 1356.7368
$ws.Range("J132").Value = 3099.9167
$ws.Range("K132").Value = 4070.2104
$ws.Range("L132").Value = 9299.750100000001
$ws.Range("M132").Value = -1540.2104
$ws.Range("N132").Value = -14359.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1089.3948
$ws.Range("I134").Value = 963.4783
$ws.Range("J134").Value = 2330.5715
$ws.Range("K134").Value = 2890.4349
$ws.Range("L134").Value = 6991.7145
$ws.Range("M134").Value = -355.4349000000002
$ws.Range("N134").Value = -12061.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 7686
$ws.Range("I53").Value = 8539
$ws.Range("K53").Value = 8539
$ws.Range("M53").Value = -7908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4296.4614
$ws.Range("I70").Value = 4089.3333
$ws.Range("K70").Value = 4089.3333
$ws.Range("M70").Value = -3819.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4296.4614
$ws.Range("I73").Value = 4089.3333
$ws.Range("K73").Value = 4089.3333
$ws.Range("M73").Value = -3153.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1996.525
$ws.Range("I132").Value = 1792.7576
$ws.Range("J132").Value = 2957.1428
$ws.Range("K132").Value = 5378.2728
$ws.Range("L132").Value = 8871.428400000001
$ws.Range("M132").Value = -2848.2728
$ws.Range("N132").Value = -13931.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 10020
$ws.Range("I34").Value = 10020
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 10020
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -9848
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2006.0667
$ws.Range("J46").Value = 2457.1428
$ws.Range("L46").Value = 2457.1428
$ws.Range("N46").Value = -2833.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 34255.5
$ws.Range("J127").Value = 34255.5
$ws.Range("L127").Value = 34255.5
$ws.Range("N127").Value = -44175.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1539.9667
$ws.Range("I132").Value = 1379.3962
$ws.Range("J132").Value = 2755.7144
$ws.Range("K132").Value = 4138.188599999999
$ws.Range("L132").Value = 8267.143199999999
$ws.Range("M132").Value = -1608.188599999999
$ws.Range("N132").Value = -13327.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2102.3901
$ws.Range("I136").Value = 1202.579
$ws.Range("J136").Value = 13500
$ws.Range("K136").Value = 3607.737
$ws.Range("L136").Value = 40500
$ws.Range("M136").Value = -1057.737
$ws.Range("N136").Value = -45600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10253.857
$ws.Range("J41").Value = 10253.857
$ws.Range("L41").Value = 10253.857
$ws.Range("N41").Value = -11033.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 618.0833
$ws.Range("I81").Value = 657
$ws.Range("J81").Value = 190
$ws.Range("K81").Value = 1314
$ws.Range("L81").Value = 380
$ws.Range("M81").Value = -253
$ws.Range("N81").Value = -2502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 618.0833
$ws.Range("I84").Value = 657
$ws.Range("J84").Value = 190
$ws.Range("K84").Value = 6570
$ws.Range("L84").Value = 1900
$ws.Range("M84").Value = -1266
$ws.Range("N84").Value = -12508

